$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue "D2" "279.08"
Set-TextValue "E2" "0.48%"
Set-TextValue "E3" "0.62%"
Set-TextValue "D4" "4.847"
Set-TextValue "E4" "-0.82%"
Set-TextValue "D5" "0.06372"
Set-TextValue "E5" "0.14%"
Set-TextValue "D6" "7.028"
Set-TextValue "E6" "0.89%"
Set-TextValue "D7" "1.304"
Set-TextValue "E7" "1.60%"
Set-TextValue "D8" "0.8952"
Set-TextValue "E8" "1.56%"
Set-TextValue "D9" "0.1534"
Set-TextValue "E9" "0.99%"
Set-TextValue "D10" "0.05808"
Set-TextValue "E10" "13.97%"
Set-TextValue "E11" "-1.41%"
Set-TextValue "E12" "-2.69%"
Set-TextValue "D13" "0.08996"
Set-TextValue "E13" "-0.24%"
Set-TextValue "D14" "0.001563"
Set-TextValue "E14" "0.03%"
Set-TextValue "D15" "0.0006391"
Set-TextValue "E15" "-0.20%"
Set-TextValue "D16" "0.006054"
Set-TextValue "E16" "0.45%"
Set-TextValue "E17" "0.60%"
Set-TextValue "D18" "3.327"
Set-TextValue "E18" "0.48%"
Set-TextValue "D19" "2.230"
Set-TextValue "E19" "-1.83%"
Set-TextValue "E21" "1.07%"
Set-TextValue "D22" "3.899"
Set-TextValue "E22" "0.00%"
Set-TextValue "D23" "0.1504"
Set-TextValue "E23" "8.96%"
Set-TextValue "D24" "0.04400"
Set-TextValue "D25" "0.001175"
Set-TextValue "E25" "0.49%"
Set-TextValue "D26" "0.004280"
Set-TextValue "E26" "10.66%"
Set-TextValue "D28" "0.0001180"
Set-TextValue "E28" "-1.67%"
Set-TextValue "D29" "0.0001654"
Set-TextValue "E29" "-14.57%"
Set-TextValue "D40" "0.04074"
Set-TextValue "E40" "-1.78%"
Set-TextValue "D41" "0.006615"
Set-TextValue "E41" "-3.39%"
Set-TextValue "D42" "0.1394"
Set-TextValue "E42" "18.30%"
Set-TextValue "D43" "0.002110"
Set-TextValue "E43" "4.46%"
Set-TextValue "D44" "0.01095"
Set-TextValue "E44" "-2.43%"
Set-TextValue "D45" "0.00005533"
Set-TextValue "E45" "6.88%"
Set-TextValue "D46" "1.628"
Set-TextValue "E46" "9.53%"
Set-TextValue "D47" "0.01850"
Set-TextValue "E47" "-8.64%"
